# T460 - update trip-sheet (foaie de parcurs) daily km / destination data
# for Alex Bora, B 151 VGT, mai 2022.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Km initiali (starting odometer reading)
$ws.Range("B12").Value = 137862

# Day 2
$ws.Range("B15").Value = 85
$ws.Range("C15").Value = "Cluj-Apahida"
$ws.Range("D15").Value = "Interes Serviciu"

# Day 3
$ws.Range("B16").Value = 47
$ws.Range("C16").Value = "Cluj-Cluj"
$ws.Range("D16").Value = "Interes Serviciu"

# Day 4
$ws.Range("B17").Value = 92
$ws.Range("C17").Value = "Cluj-Bontida"
$ws.Range("D17").Value = "Interes Serviciu"

# Day 6
$ws.Range("B19").Value = 356
$ws.Range("C19").Value = "Cluj-Baia-Mare"
$ws.Range("D19").Value = "Interes Serviciu"

# Day 9
$ws.Range("B22").Value = 421
$ws.Range("C22").Value = "Cluj-Satu-Mare"
$ws.Range("D22").Value = "Interes Serviciu"

# Day 10
$ws.Range("B23").Value = 121
$ws.Range("C23").Value = "Cluj-Turda"
$ws.Range("D23").Value = "Interes Serviciu"

# Day 11
$ws.Range("B24").Value = 92
$ws.Range("C24").Value = "Cluj-Bontida"
$ws.Range("D24").Value = "Interes Serviciu"

# Day 12
$ws.Range("B25").Value = 30
$ws.Range("C25").Value = "Acasa-Birou"
$ws.Range("D25").Value = " "

# Day 13
$ws.Range("B26").Value = 101
$ws.Range("C26").Value = "Cluj-Dej"
$ws.Range("D26").Value = "Interes Serviciu"

# Day 16
$ws.Range("B29").Value = 421
$ws.Range("C29").Value = "Cluj-Satu-Mare"
$ws.Range("D29").Value = "Interes Serviciu"

# Day 17
$ws.Range("B30").Value = 30
$ws.Range("C30").Value = "Acasa-Birou"
$ws.Range("D30").Value = " "

# Day 18
$ws.Range("B31").Value = 30
$ws.Range("C31").Value = "Acasa-Birou"
$ws.Range("D31").Value = " "

# Day 19
$ws.Range("B32").Value = 30
$ws.Range("C32").Value = "Acasa-Birou"
$ws.Range("D32").Value = " "

# Day 20
$ws.Range("B33").Value = 47
$ws.Range("C33").Value = "Cluj-Cluj"
$ws.Range("D33").Value = "Interes Serviciu"

# Day 23
$ws.Range("B36").Value = 30
$ws.Range("C36").Value = "Acasa-Birou"
$ws.Range("D36").Value = " "

# Day 24
$ws.Range("B37").Value = 30
$ws.Range("C37").Value = "Acasa-Birou"
$ws.Range("D37").Value = " "

# Day 25
$ws.Range("B38").Value = 421
$ws.Range("C38").Value = "Cluj-Satu-Mare"
$ws.Range("D38").Value = "Interes Serviciu"

# Day 26
$ws.Range("B39").Value = 30
$ws.Range("C39").Value = "Acasa-Birou"
$ws.Range("D39").Value = " "

# Day 27
$ws.Range("B40").Value = 85
$ws.Range("C40").Value = "Cluj-Apahida"
$ws.Range("D40").Value = "Interes Serviciu"

# Day 30
$ws.Range("B43").Value = 92
$ws.Range("C43").Value = "Cluj-Bontida"
$ws.Range("D43").Value = "Interes Serviciu"

# Day 31
$ws.Range("B44").Value = 356
$ws.Range("C44").Value = "Cluj-Baia-Mare"
$ws.Range("D44").Value = "Interes Serviciu"

# Totals: Km parcursi (sum of daily km) and final odometer reading
$ws.Range("B45").Value = 2947
$ws.Range("B46").Value = 140809
